$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K1").Value = "GGZRLE"
$ws.Range("N1").Value = "Min"

$ws.Range("K2").Value = 12802
$ws.Range("N2").Formula = "=SMALL(B2:K2,1)"
$ws.Range("K3").Value = 266946
$ws.Range("N3").Formula = "=SMALL(B3:K3,1)"
$ws.Range("K4").Value = 8932
$ws.Range("N4").Formula = "=SMALL(B4:K4,1)"
$ws.Range("K5").Value = 167429
$ws.Range("N5").Formula = "=SMALL(B5:K5,1)"
$ws.Range("K6").Value = 41241
$ws.Range("N6").Formula = "=SMALL(B6:K6,1)"
$ws.Range("K7").Value = 173420
$ws.Range("N7").Formula = "=SMALL(B7:K7,1)"
$ws.Range("K8").Value = 393296
$ws.Range("N8").Formula = "=SMALL(B8:K8,1)"
$ws.Range("K9").Value = 14527
$ws.Range("N9").Formula = "=SMALL(B9:K9,1)"
$ws.Range("K10").Value = 127720
$ws.Range("N10").Formula = "=SMALL(B10:K10,1)"
$ws.Range("K11").Value = 179129
$ws.Range("N11").Formula = "=SMALL(B11:K11,1)"
$ws.Range("K12").Value = 77909
$ws.Range("N12").Formula = "=SMALL(B12:K12,1)"
$ws.Range("K13").Value = 108878
$ws.Range("N13").Formula = "=SMALL(B13:K13,1)"
$ws.Range("K14").Value = 129331
$ws.Range("N14").Formula = "=SMALL(B14:K14,1)"
$ws.Range("K15").Value = 22786
$ws.Range("N15").Formula = "=SMALL(B15:K15,1)"
$ws.Range("K16").Value = 153400
$ws.Range("N16").Formula = "=SMALL(B16:K16,1)"
$ws.Range("K17").Value = 210992
$ws.Range("N17").Formula = "=SMALL(B17:K17,1)"
$ws.Range("K18").Value = 250200
$ws.Range("N18").Formula = "=SMALL(B18:K18,1)"
$ws.Range("K19").Value = 196090
$ws.Range("N19").Formula = "=SMALL(B19:K19,1)"
$ws.Range("K20").Value = 784745
$ws.Range("N20").Formula = "=SMALL(B20:K20,1)"
$ws.Range("K21").Value = 65428
$ws.Range("N21").Formula = "=SMALL(B21:K21,1)"
$ws.Range("K22").Value = 231366
$ws.Range("N22").Formula = "=SMALL(B22:K22,1)"
$ws.Range("K23").Value = 62966
$ws.Range("N23").Formula = "=SMALL(B23:K23,1)"
$ws.Range("K24").Value = 244824
$ws.Range("N24").Formula = "=SMALL(B24:K24,1)"
$ws.Range("K25").Value = 458717
$ws.Range("N25").Formula = "=SMALL(B25:K25,1)"
$ws.Range("K26").Value = 100494
$ws.Range("N26").Formula = "=SMALL(B26:K26,1)"
$ws.Range("K27").Value = 3411
$ws.Range("N27").Formula = "=SMALL(B27:K27,1)"
$ws.Range("K28").Value = 11542
$ws.Range("N28").Formula = "=SMALL(B28:K28,1)"
$ws.Range("K29").Value = 3740
$ws.Range("N29").Formula = "=SMALL(B29:K29,1)"
$ws.Range("K30").Value = 639672
$ws.Range("N30").Formula = "=SMALL(B30:K30,1)"
$ws.Range("K31").Value = 251310
$ws.Range("N31").Formula = "=SMALL(B31:K31,1)"
$ws.Range("K32").Value = 181679
$ws.Range("N32").Formula = "=SMALL(B32:K32,1)"
$ws.Range("K33").Value = 152334
$ws.Range("N33").Formula = "=SMALL(B33:K33,1)"
$ws.Range("K34").Value = 14874
$ws.Range("N34").Formula = "=SMALL(B34:K34,1)"
$ws.Range("K35").Value = 101893
$ws.Range("N35").Formula = "=SMALL(B35:K35,1)"
$ws.Range("K36").Value = 23025
$ws.Range("N36").Formula = "=SMALL(B36:K36,1)"
$ws.Range("K37").Value = 35822
$ws.Range("N37").Formula = "=SMALL(B37:K37,1)"
$ws.Range("K38").Value = 159881
$ws.Range("N38").Formula = "=SMALL(B38:K38,1)"
$ws.Range("K39").Value = 151788
$ws.Range("N39").Formula = "=SMALL(B39:K39,1)"
$ws.Range("K40").Value = 65902
$ws.Range("N40").Formula = "=SMALL(B40:K40,1)"
$ws.Range("K41").Value = 240145
$ws.Range("N41").Formula = "=SMALL(B41:K41,1)"
$ws.Range("K42").Value = 1434082
$ws.Range("N42").Formula = "=SMALL(B42:K42,1)"
$ws.Range("K43").Value = 17735
$ws.Range("N43").Formula = "=SMALL(B43:K43,1)"
$ws.Range("K44").Value = 22717
$ws.Range("N44").Formula = "=SMALL(B44:K44,1)"
$ws.Range("K45").Value = 15641
$ws.Range("N45").Formula = "=SMALL(B45:K45,1)"
$ws.Range("K46").Value = 9186
$ws.Range("N46").Formula = "=SMALL(B46:K46,1)"
$ws.Range("K47").Value = 124575
$ws.Range("N47").Formula = "=SMALL(B47:K47,1)"
$ws.Range("K48").Value = 10571
$ws.Range("N48").Formula = "=SMALL(B48:K48,1)"
$ws.Range("K49").Value = 27381
$ws.Range("N49").Formula = "=SMALL(B49:K49,1)"
$ws.Range("K50").Value = 270421
$ws.Range("N50").Formula = "=SMALL(B50:K50,1)"
$ws.Range("K51").Value = 237315
$ws.Range("N51").Formula = "=SMALL(B51:K51,1)"
$ws.Range("K52").Value = 188996
$ws.Range("N52").Formula = "=SMALL(B52:K52,1)"
